$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1 (RF 17 paragraph): "Para isso, basta informar qual palavra está
# faltando no dicionário. Assim que fizer..." becomes "Para isso, é
# obrigatório informar o nome da palavra que se deseja cadastrar. Será
# solicitado também uma descrição do sentido denotativo da palavra. Assim
# que fizer..."
# ---------------------------------------------------------------------
$old1 = "Para isso, basta informar qual palavra está faltando no dicionário. Assim que fizer a solicitação, esta irá para uma área administrativa na qual passará por revisão."
$new1 = "Para isso, é obrigatório informar o nome da palavra que se deseja cadastrar. Será solicitado também uma descrição do sentido denotativo da palavra. Assim que fizer a solicitação, esta irá para uma área administrativa na qual passará por revisão."

$rng1 = $d.Content
$rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------
# Edit 2 (RF 18 paragraph): "Segue o mesmo processo de solicitação de
# edição de palavras no dicionário (RF 17)." becomes "... (RF 17), porém
# a inclusão de descrição é opcional." -- and the _GoBack bookmark, which
# used to sit at the very end of the paragraph, ends up right before the
# final full stop instead of after it.
# ---------------------------------------------------------------------
$old2 = "Segue o mesmo processo de solicitação de edição de palavras no dicionário (RF 17)."
$new2 = "Segue o mesmo processo de solicitação de edição de palavras no dicionário (RF 17), porém a inclusão de descrição é opcional."

$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

$bm = $d.Bookmarks("_GoBack")
$bmPos = $bm.Start
$bm.Delete()
$newBmRange = $d.Range($bmPos - 1, $bmPos - 1)
$d.Bookmarks.Add("_GoBack", $newBmRange) | Out-Null

Write-Output "edits applied"
